$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Over Due" and the columns
# after it all shift one place to the right: N->O, O->P, P->Q).
$ws.Columns("N:N").Insert()

# The freshly inserted column takes on the width of its left neighbour
# (column M), just like Excel does when a column is inserted.
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet with the cursor on S9
# (it was "Transactions" with the cursor on E6/E7 before).
$ws.Activate()
$ws.Range("S9").Select()
